# mySQL_Create/types.xlsx — "Added UNIQUE constraints to db definition,
# changed classification_definition."
#
# The sheet documents the columns of a "types" table. The `name` column's
# row (row 3) gets an extra "Other 1" constraint value of UNIQUE, which
# previously was blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the UNIQUE constraint on the "name" row (row 3, column E = "Other 1") ---
# Give E3 the same look (font/border/fill) as its row-neighbour C3 (the
# "Datatype" cell for that row) before writing the new text into it, since
# both are plain unbordered body cells in that row.
$ws.Range("E3").Value = "UNIQUE"
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Selection moved from the old E12 (out of range) to E4 ---
$null = $ws.Range("E4").Select()
